# Reorder the Parameter column (A) and set corresponding Value column (B)
# per the target state captured from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("param_TimeStep_starting_index", "25"),
  @("param_demand1_op_cost_starting_index", "0"),
  @("param_demand1_inv_cost_starting_index", "0"),
  @("param_demand2_inv_cost_starting_index", "0"),
  @("param_demand2_op_cost_starting_index", "0"),
  @("param_Q_net1_demand2_starting_index", "0"),
  @("param_net1_sell_thermal_starting_index", "0"),
  @("param_net1_buy_electric_starting_index", "19.2"),
  @("param_net1_sell_electric_starting_index", "0"),
  @("param_net1_emissions_starting_index", "22.4"),
  @("param_P_net1_bat2_starting_index", "0"),
  @("param_P_net1_heat_pump2_starting_index", "20"),
  @("param_P_net1_bat1_starting_index", "0"),
  @("param_P_net1_charging_station1_starting_index", "0"),
  @("param_P_to_net1_starting_index", "0"),
  @("param_P_net1_demand2_starting_index", "0"),
  @("param_P_from_net1_starting_index", "40"),
  @("param_P_net1_demand1_starting_index", "0"),
  @("param_Q_from_net1_starting_index", "0"),
  @("param_Q_to_net1_starting_index", "0"),
  @("param_net1_inv_cost_starting_index", "0"),
  @("param_P_net1_charging_station2_starting_index", "0"),
  @("param_P_net1_heat_pump1_starting_index", "20"),
  @("param_Q_net1_demand1_starting_index", "0"),
  @("param_net1_buy_thermal_starting_index", "0"),
  @("param_net2_buy_electric_starting_index", "248.8595520199684"),
  @("param_P_net2_bat1_starting_index", "0"),
  @("param_net2_inv_cost_starting_index", "0"),
  @("param_net2_buy_thermal_starting_index", "200.9926640149763"),
  @("param_net2_sell_thermal_starting_index", "0"),
  @("param_P_net2_bat2_starting_index", "0"),
  @("param_P_net2_charging_station1_starting_index", "0"),
  @("param_Q_from_net2_starting_index", "1004.963320074881"),
  @("param_P_net2_heat_pump2_starting_index", "0"),
  @("param_P_from_net2_starting_index", "622.1488800499209"),
  @("param_P_net2_demand2_starting_index", "500"),
  @("param_Q_net2_demand2_starting_index", "757.24"),
  @("param_net2_sell_electric_starting_index", "0"),
  @("param_P_to_net2_starting_index", "0"),
  @("param_P_net2_charging_station2_starting_index", "0"),
  @("param_Q_to_net2_starting_index", "0"),
  @("param_net2_emissions_starting_index", "589.5945696459272"),
  @("param_P_net2_heat_pump1_starting_index", "0"),
  @("param_P_net2_demand1_starting_index", "122.1488800499209"),
  @("param_Q_net2_demand1_starting_index", "247.7233200748814"),
  @("param_P_pv1_demand2_starting_index", "0"),
  @("param_pv1_op_cost_starting_index", "1"),
  @("param_pv1_emissions_starting_index", "0"),
  @("param_P_pv1_demand1_starting_index", "0"),
  @("param_P_pv1_bat1_starting_index", "0"),
  @("param_P_pv1_charging_station2_starting_index", "0"),
  @("param_P_pv1_bat2_starting_index", "0"),
  @("param_P_pv1_net2_starting_index", "0"),
  @("param_P_pv1_charging_station1_starting_index", "0"),
  @("param_P_from_pv1_starting_index", "0"),
  @("param_P_pv1_heat_pump2_starting_index", "0"),
  @("param_P_pv1_heat_pump1_starting_index", "0"),
  @("param_P_pv1_net1_starting_index", "0"),
  @("param_pv1_inv_cost_starting_index", "0"),
  @("param_P_pv2_charging_station1_starting_index", "0"),
  @("param_P_pv2_net1_starting_index", "0"),
  @("param_P_from_pv2_starting_index", "1.8"),
  @("param_P_pv2_bat1_starting_index", "0"),
  @("param_P_pv2_bat2_starting_index", "0"),
  @("param_P_pv2_demand2_starting_index", "0"),
  @("param_P_pv2_heat_pump1_starting_index", "0"),
  @("param_P_pv2_heat_pump2_starting_index", "0"),
  @("param_pv2_inv_cost_starting_index", "0"),
  @("param_P_pv2_charging_station2_starting_index", "0"),
  @("param_P_pv2_demand1_starting_index", "1.8"),
  @("param_P_pv2_net2_starting_index", "0"),
  @("param_pv2_op_cost_starting_index", "1"),
  @("param_pv2_emissions_starting_index", "0.8999999999999999"),
  @("param_bat1_K_ch_starting_index", "0"),
  @("param_P_bat1_net1_starting_index", "0"),
  @("param_P_bat1_heat_pump1_starting_index", "0"),
  @("param_P_bat1_demand2_starting_index", "0"),
  @("param_bat1_cumulated_aging_starting_index", "3.166666666666667E-06"),
  @("param_bat1_op_cost_starting_index", "1"),
  @("param_bat1_emissions_starting_index", "0"),
  @("param_bat1_K_dis_starting_index", "1"),
  @("param_bat1_integer_starting_index", "-0"),
  @("param_bat1_SOC_starting_index", "0.3"),
  @("param_P_bat1_charging_station1_starting_index", "0"),
  @("param_P_bat1_heat_pump2_starting_index", "0"),
  @("param_P_to_bat1_starting_index", "0"),
  @("param_bat1_inv_cost_starting_index", "0"),
  @("param_P_from_bat1_starting_index", "0"),
  @("param_P_bat1_demand1_starting_index", "0"),
  @("param_P_bat1_charging_station2_starting_index", "0"),
  @("param_P_bat1_net2_starting_index", "0"),
  @("param_bat1_SOC_max_starting_index", "0.9999968333333333"),
  @("param_P_bat2_net2_starting_index", "0"),
  @("param_P_bat2_demand1_starting_index", "0"),
  @("param_bat2_SOC_starting_index", "0.5"),
  @("param_P_bat2_charging_station1_starting_index", "0"),
  @("param_bat2_op_cost_starting_index", "1"),
  @("param_P_bat2_net1_starting_index", "0"),
  @("param_bat2_inv_cost_starting_index", "0"),
  @("param_bat2_emissions_starting_index", "0"),
  @("param_bat2_cumulated_aging_starting_index", "0"),
  @("param_P_bat2_demand2_starting_index", "0"),
  @("param_P_to_bat2_starting_index", "0"),
  @("param_P_bat2_charging_station2_starting_index", "0"),
  @("param_P_bat2_heat_pump1_starting_index", "0"),
  @("param_P_bat2_heat_pump2_starting_index", "0"),
  @("param_P_from_bat2_starting_index", "0"),
  @("param_bat2_K_dis_starting_index", "1"),
  @("param_bat2_K_ch_starting_index", "0"),
  @("param_bat2_SOC_max_starting_index", "1"),
  @("param_bat2_integer_starting_index", "0"),
  @("param_Q_CHP1_demand1_starting_index", "0"),
  @("param_P_CHP1_bat2_starting_index", "0"),
  @("param_P_CHP1_charging_station2_starting_index", "0"),
  @("param_P_CHP1_demand2_starting_index", "0"),
  @("param_CHP1_emissions_starting_index", "4.83"),
  @("param_P_CHP1_heat_pump1_starting_index", "0"),
  @("param_P_CHP1_bat1_starting_index", "0"),
  @("param_P_CHP1_net2_starting_index", "0"),
  @("param_CHP1_inv_cost_starting_index", "0"),
  @("param_P_CHP1_net1_starting_index", "0"),
  @("param_Q_CHP1_net2_starting_index", "0"),
  @("param_P_CHP1_charging_station1_starting_index", "0"),
  @("param_P_CHP1_heat_pump2_starting_index", "0"),
  @("param_P_from_CHP1_starting_index", "20"),
  @("param_Q_CHP1_net1_starting_index", "0"),
  @("param_P_CHP1_demand1_starting_index", "20"),
  @("param_CHP1_fuel_cons_starting_index", "2.1"),
  @("param_CHP1_op_cost_starting_index", "10.5"),
  @("param_Q_from_CHP1_starting_index", "40"),
  @("param_Q_CHP1_demand2_starting_index", "40"),
  @("param_P_CHP2_net2_starting_index", "0"),
  @("param_P_CHP2_bat2_starting_index", "0"),
  @("param_P_from_CHP2_starting_index", "20"),
  @("param_CHP2_op_cost_starting_index", "10.5"),
  @("param_Q_CHP2_net1_starting_index", "0"),
  @("param_P_CHP2_heat_pump2_starting_index", "0"),
  @("param_Q_CHP2_net2_starting_index", "0"),
  @("param_Q_CHP2_demand1_starting_index", "0"),
  @("param_CHP2_fuel_cons_starting_index", "2.1"),
  @("param_P_CHP2_charging_station1_starting_index", "0"),
  @("param_CHP2_inv_cost_starting_index", "0"),
  @("param_P_CHP2_bat1_starting_index", "0"),
  @("param_P_CHP2_demand1_starting_index", "20"),
  @("param_Q_from_CHP2_starting_index", "40"),
  @("param_Q_CHP2_demand2_starting_index", "40"),
  @("param_P_CHP2_heat_pump1_starting_index", "0"),
  @("param_P_CHP2_net1_starting_index", "0"),
  @("param_CHP2_emissions_starting_index", "4.83"),
  @("param_P_CHP2_demand2_starting_index", "0"),
  @("param_P_CHP2_charging_station2_starting_index", "0"),
  @("param_Q_solar_th1_net2_starting_index", "0"),
  @("param_Q_from_solar_th1_starting_index", "0"),
  @("param_Q_solar_th1_demand1_starting_index", "0"),
  @("param_Q_solar_th1_net1_starting_index", "0"),
  @("param_solar_th1_op_cost_starting_index", "1"),
  @("param_solar_th1_inv_cost_starting_index", "0"),
  @("param_solar_th1_emissions_starting_index", "0"),
  @("param_Q_solar_th1_demand2_starting_index", "0"),
  @("param_Q_solar_th2_net1_starting_index", "0"),
  @("param_Q_solar_th2_net2_starting_index", "0"),
  @("param_Q_solar_th2_demand2_starting_index", "1.2"),
  @("param_solar_th2_emissions_starting_index", "0.6000000000000001"),
  @("param_Q_from_solar_th2_starting_index", "1.2"),
  @("param_solar_th2_op_cost_starting_index", "1"),
  @("param_Q_solar_th2_demand1_starting_index", "0"),
  @("param_solar_th2_inv_cost_starting_index", "0"),
  @("param_P_pvt1_bat2_starting_index", "0"),
  @("param_P_from_pvt1_starting_index", "0"),
  @("param_pvt1_emissions_starting_index", "0"),
  @("param_pvt1_inv_cost_starting_index", "0"),
  @("param_P_pvt1_bat1_starting_index", "0"),
  @("param_P_pvt1_net2_starting_index", "0"),
  @("param_pvt1_op_cost_starting_index", "1"),
  @("param_P_pvt1_heat_pump2_starting_index", "0"),
  @("param_P_pvt1_charging_station2_starting_index", "0"),
  @("param_Q_pvt1_net1_starting_index", "0"),
  @("param_Q_from_pvt1_starting_index", "0"),
  @("param_P_pvt1_net1_starting_index", "0"),
  @("param_Q_pvt1_demand2_starting_index", "0"),
  @("param_P_pvt1_demand1_starting_index", "0"),
  @("param_P_pvt1_heat_pump1_starting_index", "0"),
  @("param_Q_pvt1_net2_starting_index", "0"),
  @("param_P_pvt1_charging_station1_starting_index", "0"),
  @("param_Q_pvt1_demand1_starting_index", "0"),
  @("param_P_pvt1_demand2_starting_index", "0"),
  @("param_P_pvt2_bat1_starting_index", "0"),
  @("param_P_pvt2_demand2_starting_index", "0"),
  @("param_P_pvt2_net1_starting_index", "0"),
  @("param_Q_from_pvt2_starting_index", "1.56"),
  @("param_P_from_pvt2_starting_index", "1.2"),
  @("param_P_pvt2_charging_station2_starting_index", "0"),
  @("param_P_pvt2_heat_pump2_starting_index", "0"),
  @("param_P_pvt2_bat2_starting_index", "0"),
  @("param_pvt2_emissions_starting_index", "0.7800000000000001"),
  @("param_Q_pvt2_demand1_starting_index", "0"),
  @("param_P_pvt2_charging_station1_starting_index", "0"),
  @("param_P_pvt2_demand1_starting_index", "1.2"),
  @("param_Q_pvt2_net2_starting_index", "0"),
  @("param_P_pvt2_net2_starting_index", "0"),
  @("param_P_pvt2_heat_pump1_starting_index", "0"),
  @("param_Q_pvt2_net1_starting_index", "0"),
  @("param_pvt2_inv_cost_starting_index", "0"),
  @("param_Q_pvt2_demand2_starting_index", "1.56"),
  @("param_pvt2_op_cost_starting_index", "1"),
  @("param_charging_station1_inv_cost_starting_index", "0"),
  @("param_charging_station1_op_cost_starting_index", "0"),
  @("param_charging_station1_emissions_starting_index", "0"),
  @("param_charging_station2_emissions_starting_index", "0"),
  @("param_charging_station2_inv_cost_starting_index", "0"),
  @("param_charging_station2_op_cost_starting_index", "0"),
  @("param_Q_heat_pump1_net1_starting_index", "0"),
  @("param_heat_pump1_op_cost_starting_index", "8.561643835616438"),
  @("param_heat_pump1_emissions_starting_index", "2.76"),
  @("param_Q_from_heat_pump1_starting_index", "80"),
  @("param_Q_to_heat_pump1_starting_index", "0"),
  @("param_P_from_heat_pump1_starting_index", "0"),
  @("param_heat_pump1_inv_cost_starting_index", "0"),
  @("param_Q_heat_pump1_net2_starting_index", "0"),
  @("param_Q_heat_pump1_demand1_starting_index", "0"),
  @("param_Q_heat_pump1_demand2_starting_index", "80"),
  @("param_P_to_heat_pump1_starting_index", "20"),
  @("param_P_to_heat_pump2_starting_index", "20"),
  @("param_Q_from_heat_pump2_starting_index", "80"),
  @("param_heat_pump2_emissions_starting_index", "2.76"),
  @("param_P_from_heat_pump2_starting_index", "0"),
  @("param_Q_heat_pump2_demand2_starting_index", "80"),
  @("param_Q_heat_pump2_net1_starting_index", "0"),
  @("param_Q_heat_pump2_net2_starting_index", "0"),
  @("param_Q_heat_pump2_demand1_starting_index", "0"),
  @("param_heat_pump2_op_cost_starting_index", "8.561643835616438"),
  @("param_heat_pump2_inv_cost_starting_index", "0"),
  @("param_Q_to_heat_pump2_starting_index", "0"),
  @("param_total_emissions_starting_index", "628.0745696459272"),
  @("param_total_sell_starting_index", "0"),
  @("param_total_buy_starting_index", "469.0522160349446"),
  @("param_total_operation_cost_starting_index", "469.0522160349446")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $rows[$i][0]
  $ws.Cells.Item($r, 2).Value = [double]$rows[$i][1]
}